$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30
$prev = $row - 1

$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"

$ws.Cells.Item($row, 6).Value = "Mumbai City"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Punjab"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.33
$ws.Cells.Item($row, 11).Value = "01/11/2023 03:12"
$ws.Cells.Item($row, 12).Value = 1.25
$ws.Cells.Item($row, 13).Value = "02/11/2023 15:02"
$ws.Cells.Item($row, 14).Value = 5.26
$ws.Cells.Item($row, 15).Value = "01/11/2023 03:12"
$ws.Cells.Item($row, 16).Value = 6.17
$ws.Cells.Item($row, 17).Value = "02/11/2023 15:02"
$ws.Cells.Item($row, 18).Value = 8.56
$ws.Cells.Item($row, 19).Value = "01/11/2023 03:12"
$ws.Cells.Item($row, 20).Value = 10.51
$ws.Cells.Item($row, 21).Value = "02/11/2023 15:02"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/mumbai-city-minerva-punjab/GxhgvCoc/"

# A30 (Indice) and E30 (data_partida) carry special cell formatting (s="1" and s="2")
# matching the preceding row, so copy the format over before setting their values.
$ws.Range($ws.Cells.Item($prev, 1), $ws.Cells.Item($prev, 1)).Copy()
$ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 1)).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 29

$ws.Range($ws.Cells.Item($prev, 5), $ws.Cells.Item($prev, 5)).Copy()
$ws.Range($ws.Cells.Item($row, 5), $ws.Cells.Item($row, 5)).PasteSpecial(-4122)
$ws.Cells.Item($row, 5).Value = 45232.64583333334

$excel.CutCopyMode = 0
